# Snack-Automat "Test cases.xlsx" - finish test 2
#
# Sheet2 holds two small tables:
#   Table5 (B16:E21)  -> Nr. | Test | Description | ..      (4th column header was a placeholder)
#   Table6 (I16:J22)  -> Test case number | 1                (2nd column header was test-case "1")
#
# Test case "1" is done; this pass finishes test case "2":
#   - Table5's placeholder header becomes "Status" and every test row is marked "OK"
#   - Table6's header becomes "2" (the table now documents test case 2's number)
#   - the view is scrolled down a bit and the selection left on the Table6 area

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# --- Table5 (Nr./Test/Description/Status) ------------------------------
# Renaming the header cell also renames the underlying ListObject column.
$ws.Range("E16").Value = "Status"

# Mark every existing test-case row as finished/OK.
$ws.Range("E17").Value = "OK"
$ws.Range("E18").Value = "OK"
$ws.Range("E19").Value = "OK"
$ws.Range("E20").Value = "OK"
$ws.Range("E21").Value = "OK"

# The new "Status" column only needs to be narrow.
$ws.Columns.Item(5).ColumnWidth = 9.5

# --- Table6 (Test case number / 2) --------------------------------------
$ws.Range("J16").Value = "2"

# Leave the view scrolled to where work left off.
$ws.Range("J25").Select()
